$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.582.02"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "2.065.56"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'242.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").Value = "'0.661"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'52.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.59%  "
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "'0.358"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.20%  "
$ws.Range("D11").Value = "'0.0750"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.60%  "
$ws.Range("D13").Value = "'0.900"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "'14.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.60%  "
$ws.Range("D15").Value = "2.364.61"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "'5.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.25%  "
$ws.Range("D17").Value = "2.047.58"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "36.514.16"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").Value = "'16.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -16.16%  "
$ws.Range("D20").Value = "'71.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.80%  "
$ws.Range("D21").Value = "0.0₃0862"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.55%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'235.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  -5.36%  "
$ws.Range("D26").Value = "'9.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.41%  "
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("D28").Value = "'163.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.00%  "
$ws.Range("D29").Value = "'20.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").Value = "'5.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'1.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").Value = "'0.0596"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.80%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'2.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "'0.0821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.79%  "
$ws.Range("D39").Value = "'1.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.47%  "
$ws.Range("E40").Value = "  -5.22%  "
$ws.Range("D41").Value = "'4.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.70%  "
$ws.Range("E42").Value = "  -3.93%  "
$ws.Range("D43").Value = "'1.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("D44").Value = "'0.0933"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.92%  "
$ws.Range("D45").Value = "'93.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.73%  "
$ws.Range("D46").Value = "1.404.24"
$ws.Range("E46").Value = "  +9.21%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'15.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.16%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.59%  "
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").Value = "2.252.93"
$ws.Range("E51").Value = "  +0.16%  "
